# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect newer counts scraped at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 96
$ws1.Range("F4").Value  = 611
$ws1.Range("F6").Value  = 9202
$ws1.Range("F9").Value  = 1178
$ws1.Range("F10").Value = 1065
$ws1.Range("F12").Value = 55
$ws1.Range("F14").Value = 249
$ws1.Range("F15").Value = 364
$ws1.Range("F18").Value = 1200

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 96
$ws4.Range("F6").Value  = 611
$ws4.Range("F8").Value  = 9202
$ws4.Range("F11").Value = 1178
$ws4.Range("F12").Value = 1065
$ws4.Range("F14").Value = 55
$ws4.Range("F16").Value = 249
$ws4.Range("F17").Value = 364
$ws4.Range("F20").Value = 1200
